# Add a new "2022-Q1" sheet (fund-level holdings) positioned right before
# the "总计" (totals) sheet, and add a corresponding summary row to "总计".

$wb = $excel.ActiveWorkbook

# --- locate the reference sheets --------------------------------------------------
$q4sheet    = $wb.Worksheets.Item("2021-Q4")
$totalSheet0 = $wb.Worksheets.Item("总计")

# --- create + place the new sheet --------------------------------------------------
$tmp = $wb.Worksheets.Add()
$tmp.Name = "2022-Q1"
$tmp.Move($null, $totalSheet0)   # move to just after 2021-Q4 / before 总计

# NOTE: sheet handles captured via .Item(...)/.Add() go stale (silently point at
# the wrong sheet) once the tab order changes underneath them (Add/Move/Delete).
# Re-resolve every handle we still need to write through *after* all of the
# structural changes (add/move) above are finished.
$newSheet   = $wb.Worksheets.Item("2022-Q1")
$totalSheet = $wb.Worksheets.Item("总计")
$q4sheet    = $wb.Worksheets.Item("2021-Q4")

# Copy the header-cell style (bold + border + centered) from an existing sheet's
# B1 cell so the new sheet's header row / index column match the workbook look.
# NOTE: some operations (ClearFormats, etc.) silently drop the Office clipboard,
# so the source cell is re-Copy()-ed right before every PasteSpecial below
# instead of relying on one Copy() lasting the whole script.

# --- header row ----------------------------------------------------------------
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = 2 + $i   # starts at column B
    $cell = $newSheet.Cells.Item(1, $col)
    $q4sheet.Range("B1").Copy()
    $cell.PasteSpecial(-4122)   # xlPasteFormats - copy style only
    $cell.Value = $headers[$i]
}

# --- data rows -------------------------------------------------------------------
# Columns B-G hold text (fund code keeps leading zeros, decimals keep their
# original formatting like "0.0260" / "0.00"); only a zero holding value (G)
# and the rank column (H) are real numbers - matches the source data.
$rows = @(
    @("007107", "太平 MSCI 香港价值增强指数A",     "1.05", "93.78", "2.48", "0.0260", $true,  8),
    @("004532", "民生加银中证港股通高股息精选指数A", "0.26", "94.88", "3.68", "0.0096", $true,  6),
    @("004533", "民生加银中证港股通高股息精选指数C", "0.10", "94.88", "3.68", "0.0037", $true,  6),
    @("007108", "太平 MSCI 香港价值增强指数C",     "0.00", "93.78", "2.48", 0,        $false, 8)
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $row = 2 + $r
    $a = $newSheet.Cells.Item($row, 1)
    $a.PasteSpecial(-4122)
    $a.Value = $r

    $bCell = $newSheet.Cells.Item($row, 2)
    $bCell.NumberFormat = "@"
    $bCell.Value = $rows[$r][0]

    $newSheet.Cells.Item($row, 3).Value = $rows[$r][1]

    $dCell = $newSheet.Cells.Item($row, 4)
    $dCell.NumberFormat = "@"
    $dCell.Value = $rows[$r][2]

    $eCell = $newSheet.Cells.Item($row, 5)
    $eCell.NumberFormat = "@"
    $eCell.Value = $rows[$r][3]

    $fCell = $newSheet.Cells.Item($row, 6)
    $fCell.NumberFormat = "@"
    $fCell.Value = $rows[$r][4]

    $gCell = $newSheet.Cells.Item($row, 7)
    if ($rows[$r][6]) {
        $gCell.NumberFormat = "@"
    }
    $gCell.Value = $rows[$r][5]

    $newSheet.Cells.Item($row, 8).Value = $rows[$r][7]
}

# --- update the "总计" (totals) sheet: insert a new top data row for 2022-Q1 -------
$totalSheet.Rows(2).Insert()
$totalSheet.Range("A2:D2").ClearFormats()

$aCell = $totalSheet.Cells.Item(2, 1)
$aCell.PasteSpecial(-4122)
$aCell.Value = 0

$totalSheet.Cells.Item(2, 2).Value = "2022-Q1"
$totalSheet.Cells.Item(2, 3).Value = 4
$totalSheet.Cells.Item(2, 4).Value = 0.04

# re-number the index column (A) for the rows that got pushed down
for ($r = 3; $r -le 7; $r++) {
    $totalSheet.Cells.Item($r, 1).Value = $r - 2
}

$excel.CutCopyMode = $false
